# Auto-generated Excel COM-interop script
# Applies scheduled-runner market data refresh to Sheets/Adamantoise_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 448.66666
$ws.Range("I55").Value = 360.44446
$ws.Range("J55").Value = 536.8889
$ws.Range("K55").Value = 360.44446
$ws.Range("L55").Value = 536.8889
$ws.Range("M55").Value = -146.44446
$ws.Range("N55").Value = -964.8889
$ws.Range("H58").Value = 3153.75
$ws.Range("J58").Value = 4166.6665
$ws.Range("L58").Value = 12499.9995
$ws.Range("N58").Value = -12799.9995
$ws.Range("H101").Value = 935.5
$ws.Range("I101").Value = 997.6
$ws.Range("J101").Value = 625
$ws.Range("K101").Value = 2992.8
$ws.Range("L101").Value = 1875
$ws.Range("M101").Value = -1370.8
$ws.Range("N101").Value = -5119
$ws.Range("H132").Value = 2010.6666
$ws.Range("I132").Value = 2095.361
$ws.Range("K132").Value = 6286.083
$ws.Range("M132").Value = -3756.083
$ws.Range("H135").Value = 1527.1052
$ws.Range("I135").Value = 1528.6111
$ws.Range("K135").Value = 13757.4999
$ws.Range("M135").Value = -11222.4999
$ws.Range("H137").Value = 1463933.6
$ws.Range("I137").Value = 1407.7391
$ws.Range("J137").Value = 3706473.5
$ws.Range("K137").Value = 4223.2173
$ws.Range("L137").Value = 11119420.5
$ws.Range("M137").Value = -1673.2173
$ws.Range("N137").Value = -11124520.5
$ws.Range("H138").Value = 1790.5
$ws.Range("I138").Value = 657.4211
$ws.Range("J138").Value = 2484.9678
$ws.Range("K138").Value = 1972.2633
$ws.Range("L138").Value = 7454.903399999999
$ws.Range("M138").Value = 3167.7367
$ws.Range("N138").Value = -17734.9034
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20375594
$ws.Range("I32").Value = 24024780
$ws.Range("J32").Value = 5498145
$ws.Range("K32").Value = 24024780
$ws.Range("L32").Value = 5498145
$ws.Range("M32").Value = -24024493
$ws.Range("N32").Value = -5498719
$ws.Range("H55").Value = 56678
$ws.Range("J55").Value = 69993
$ws.Range("L55").Value = 69993
$ws.Range("N55").Value = -70623
$ws.Range("H74").Value = 3533
$ws.Range("I74").Value = 2945.2
$ws.Range("K74").Value = 2945.2
$ws.Range("M74").Value = -2071.2
$ws.Range("H77").Value = 3533
$ws.Range("I77").Value = 2945.2
$ws.Range("K77").Value = 14726
$ws.Range("M77").Value = -10358
$ws.Range("H80").Value = 75000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 75000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H92").Value = 67996
$ws.Range("J92").Value = 67996
$ws.Range("L92").Value = 67996
$ws.Range("N92").Value = -72988
$ws.Range("H96").Value = 66380.8
$ws.Range("J96").Value = 66380.8
$ws.Range("L96").Value = 66380.8
$ws.Range("N96").Value = -71872.8
$ws.Range("H97").Value = 737.0454999999999
$ws.Range("I97").Value = 737.0454999999999
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 737.0454999999999
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -241.0454999999999
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 2882.0815
$ws.Range("I132").Value = 2605.475
$ws.Range("J132").Value = 4111.4443
$ws.Range("K132").Value = 7816.424999999999
$ws.Range("L132").Value = 12334.3329
$ws.Range("M132").Value = -5286.424999999999
$ws.Range("N132").Value = -17394.3329
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 61017.8
$ws.Range("J28").Value = 61017.8
$ws.Range("L28").Value = 61017.8
$ws.Range("N28").Value = -61507.8
$ws.Range("H31").Value = 5428.685
$ws.Range("I31").Value = 2179.36
$ws.Range("J31").Value = 8229.826999999999
$ws.Range("K31").Value = 2179.36
$ws.Range("L31").Value = 8229.826999999999
$ws.Range("M31").Value = -1884.36
$ws.Range("N31").Value = -8819.826999999999
$ws.Range("H34").Value = 5428.685
$ws.Range("I34").Value = 2179.36
$ws.Range("J34").Value = 8229.826999999999
$ws.Range("K34").Value = 2179.36
$ws.Range("L34").Value = 8229.826999999999
$ws.Range("M34").Value = -1977.36
$ws.Range("N34").Value = -8633.826999999999
$ws.Range("H69").Value = 63609.855
$ws.Range("I69").Value = 50000
$ws.Range("K69").Value = 50000
$ws.Range("M69").Value = -49251
$ws.Range("H72").Value = 63609.855
$ws.Range("I72").Value = 50000
$ws.Range("K72").Value = 150000
$ws.Range("M72").Value = -146256
$ws.Range("H132").Value = 3422.524
$ws.Range("I132").Value = 3235.5715
$ws.Range("J132").Value = 3796.4285
$ws.Range("K132").Value = 9706.7145
$ws.Range("L132").Value = 11389.2855
$ws.Range("M132").Value = -7176.7145
$ws.Range("N132").Value = -16449.2855
$ws.Range("H134").Value = 1823.6316
$ws.Range("I134").Value = 1823.6316
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5470.8948
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2935.8948
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 510.451
$ws.Range("I107").Value = 1179.1111
$ws.Range("K107").Value = 3537.3333
$ws.Range("M107").Value = -1617.3333
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1731.75
$ws.Range("I97").Value = 1410
$ws.Range("J97").Value = 2182.2
$ws.Range("K97").Value = 1410
$ws.Range("L97").Value = 2182.2
$ws.Range("M97").Value = -914
$ws.Range("N97").Value = -3174.2
$ws.Range("H107").Value = 1900
$ws.Range("I107").Value = 1900
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1900
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 20
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 21758.2
$ws.Range("I113").Value = 2126.6667
$ws.Range("K113").Value = 2126.6667
$ws.Range("M113").Value = 43.33329999999978
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1942.3846
$ws.Range("I22").Value = 1372.4445
$ws.Range("J22").Value = 3224.75
$ws.Range("K22").Value = 1372.4445
$ws.Range("L22").Value = 3224.75
$ws.Range("M22").Value = -1077.4445
$ws.Range("N22").Value = -3814.75
$ws.Range("H27").Value = 1942.3846
$ws.Range("I27").Value = 1372.4445
$ws.Range("J27").Value = 3224.75
$ws.Range("K27").Value = 1372.4445
$ws.Range("L27").Value = 3224.75
$ws.Range("M27").Value = -1265.4445
$ws.Range("N27").Value = -3438.75
$ws.Range("H46").Value = 3471.2593
$ws.Range("I46").Value = 774
$ws.Range("J46").Value = 4241.905
$ws.Range("K46").Value = 774
$ws.Range("L46").Value = 4241.905
$ws.Range("M46").Value = -586
$ws.Range("N46").Value = -4617.905
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H136").Value = 2597.6
$ws.Range("I136").Value = 2420.182
$ws.Range("J136").Value = 3085.5
$ws.Range("K136").Value = 7260.545999999999
$ws.Range("L136").Value = 9256.5
$ws.Range("M136").Value = -4710.545999999999
$ws.Range("N136").Value = -14356.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 135000
$ws.Range("J57").Value = 135000
$ws.Range("L57").Value = 135000
$ws.Range("N57").Value = -136508
$ws.Range("H62").Value = 5168
$ws.Range("J62").Value = 6500
$ws.Range("L62").Value = 6500
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 5168
$ws.Range("J65").Value = 6500
$ws.Range("L65").Value = 32500
$ws.Range("N65").Value = -38740
$ws.Range("H132").Value = 1389.7428
$ws.Range("I132").Value = 1405.8
$ws.Range("K132").Value = 4217.4
$ws.Range("M132").Value = -1687.4
